$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / value updates -------------------------------------------------

# Mã phiếu column (A2:A4) - was P2209-00001, now P2210-00001 for every row
$ws.Range("A2:A4").Value = "P2210-00001"

# Row 2: Số lô (C2) changes 11 -> 1234
$ws.Range("C2").Value = 1234

# Row 3: medicine becomes Acemuc 100mg, C3 12 -> 111, E3 100 -> 200,
# and the "Tổng giá nhập" (H3) becomes 600.000 VND
$ws.Range("B3").Value = "Acemuc 100mg"
$ws.Range("C3").Value = 111
$ws.Range("E3").Value = 200
$ws.Range("H3").Value = "600.000 VND"

# Row 4: medicine becomes Acemuc 200mg (replacing Amlordipin 5mg), C4 13 -> 222
$ws.Range("B4").Value = "Acemuc 200mg"
$ws.Range("C4").Value = 222

# --- Column width -----------------------------------------------------
# Column B's stored sheet width goes from 17 to 15 characters. Excel's
# ColumnWidth property is offset from the stored <col width> by the
# default font padding (~0.7142857 chars), so compensate here.
$ws.Columns.Item(2).ColumnWidth = 14.2857142857
